$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.21"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "25.10"
$ws.Range("D3").Style = "Normal"

$ws.Range("B4").Value = "HuobiToken"

$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.130"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "3HuobiTokenHT"

$ws.Range("B5").Value = "Cronos"

$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05753"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "4CronosCRO"

$ws.Range("B6").Value = "KuCoinToken"

$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.470"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "5KuCoinTokenKCS"

$ws.Range("B7").Value = "GateToken"

$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.119"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = "6GateTokenGT"

$ws.Range("B8").Value = "MXToken"

$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8099"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "7MXTokenMX"

$ws.Range("B9").Value = "FTXToken"

$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8404"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "8FTXTokenFTT"

$ws.Range("B10").Value = "One"

$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0005997"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "9OneONE"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1338"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06939"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03135"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.02830"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09369"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.758"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001528"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04663"
$ws.Range("D18").Style = "Normal"

$ws.Range("B19").Value = "TigerCash"

$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006087"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "BitKan"

$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001235"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "HotbitToken"

$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004281"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "NitroEx"

$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008696"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "LEO"

$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.501"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002325"
$ws.Range("D28").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03611"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006387"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1050"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002936"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007369"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005306"
$ws.Range("D45").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.2509"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

Write-Host "Applied all cell updates"
